$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, reusing the existing header style (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Per-row time_taken timestamps for rows 2-191
$timeTaken = @(
    "2021-10-05 13:41:20.813132",
    "2021-10-05 13:41:20.813143",
    "2021-10-05 13:41:20.813147",
    "2021-10-05 13:41:20.813149",
    "2021-10-05 13:41:20.813152",
    "2021-10-05 13:41:20.813155",
    "2021-10-05 13:41:20.813158",
    "2021-10-05 13:41:20.813160",
    "2021-10-05 13:41:20.813163",
    "2021-10-05 13:41:20.813166",
    "2021-10-05 13:41:20.813168",
    "2021-10-05 13:41:20.813171",
    "2021-10-05 13:41:20.813174",
    "2021-10-05 13:41:20.813176",
    "2021-10-05 13:41:20.813179",
    "2021-10-05 13:41:20.813181",
    "2021-10-05 13:41:20.813184",
    "2021-10-05 13:41:20.813187",
    "2021-10-05 13:41:20.813189",
    "2021-10-05 13:41:20.813192",
    "2021-10-05 13:41:20.813195",
    "2021-10-05 13:41:20.813197",
    "2021-10-05 13:41:20.813200",
    "2021-10-05 13:41:20.813202",
    "2021-10-05 13:41:20.813205",
    "2021-10-05 13:41:20.813208",
    "2021-10-05 13:41:20.813211",
    "2021-10-05 13:41:20.813213",
    "2021-10-05 13:41:20.813216",
    "2021-10-05 13:41:20.813219",
    "2021-10-05 13:41:20.813221",
    "2021-10-05 13:41:20.813224",
    "2021-10-05 13:41:20.813227",
    "2021-10-05 13:41:20.813230",
    "2021-10-05 13:41:20.813233",
    "2021-10-05 13:41:20.813235",
    "2021-10-05 13:41:20.813238",
    "2021-10-05 13:41:20.813240",
    "2021-10-05 13:41:20.813243",
    "2021-10-05 13:41:20.813245",
    "2021-10-05 13:41:20.813248",
    "2021-10-05 13:41:20.813251",
    "2021-10-05 13:41:20.813253",
    "2021-10-05 13:41:20.813256",
    "2021-10-05 13:41:20.813259",
    "2021-10-05 13:41:20.813262",
    "2021-10-05 13:41:20.813264",
    "2021-10-05 13:41:20.813267",
    "2021-10-05 13:41:20.813269",
    "2021-10-05 13:41:20.813272",
    "2021-10-05 13:41:20.813275",
    "2021-10-05 13:41:20.813277",
    "2021-10-05 13:41:20.813280",
    "2021-10-05 13:41:20.813282",
    "2021-10-05 13:41:20.813285",
    "2021-10-05 13:41:20.813287",
    "2021-10-05 13:41:20.813290",
    "2021-10-05 13:41:20.813292",
    "2021-10-05 13:41:20.813295",
    "2021-10-05 13:41:20.813297",
    "2021-10-05 13:41:20.813300",
    "2021-10-05 13:41:20.813302",
    "2021-10-05 13:41:20.813305",
    "2021-10-05 13:41:20.813308",
    "2021-10-05 13:41:20.813311",
    "2021-10-05 13:41:20.813314",
    "2021-10-05 13:41:20.813317",
    "2021-10-05 13:41:20.813319",
    "2021-10-05 13:41:20.813322",
    "2021-10-05 13:41:20.813324",
    "2021-10-05 13:41:20.813327",
    "2021-10-05 13:41:20.813329",
    "2021-10-05 13:41:20.813332",
    "2021-10-05 13:41:20.813334",
    "2021-10-05 13:41:20.813337",
    "2021-10-05 13:41:20.813339",
    "2021-10-05 13:41:20.813344",
    "2021-10-05 13:41:20.813347",
    "2021-10-05 13:41:20.813349",
    "2021-10-05 13:41:20.813352",
    "2021-10-05 13:41:20.813355",
    "2021-10-05 13:41:20.813357",
    "2021-10-05 13:41:20.813360",
    "2021-10-05 13:41:20.813362",
    "2021-10-05 13:41:20.813365",
    "2021-10-05 13:41:20.813367",
    "2021-10-05 13:41:20.813370",
    "2021-10-05 13:41:20.813372",
    "2021-10-05 13:41:20.813375",
    "2021-10-05 13:41:20.813378",
    "2021-10-05 13:41:20.813380",
    "2021-10-05 13:41:20.813383",
    "2021-10-05 13:41:20.813386",
    "2021-10-05 13:41:20.813389",
    "2021-10-05 13:41:20.813392",
    "2021-10-05 13:41:20.813395",
    "2021-10-05 13:41:20.813397",
    "2021-10-05 13:41:20.813400",
    "2021-10-05 13:41:20.813402",
    "2021-10-05 13:41:20.813405",
    "2021-10-05 13:41:20.813407",
    "2021-10-05 13:41:20.813410",
    "2021-10-05 13:41:20.813412",
    "2021-10-05 13:41:20.813415",
    "2021-10-05 13:41:20.813417",
    "2021-10-05 13:41:20.813420",
    "2021-10-05 13:41:20.813423",
    "2021-10-05 13:41:20.813425",
    "2021-10-05 13:41:20.813430",
    "2021-10-05 13:41:20.813433",
    "2021-10-05 13:41:20.813436",
    "2021-10-05 13:41:20.813438",
    "2021-10-05 13:41:20.813441",
    "2021-10-05 13:41:20.813443",
    "2021-10-05 13:41:20.813446",
    "2021-10-05 13:41:20.813448",
    "2021-10-05 13:41:20.813451",
    "2021-10-05 13:41:20.813453",
    "2021-10-05 13:41:20.813456",
    "2021-10-05 13:41:20.813458",
    "2021-10-05 13:41:20.813461",
    "2021-10-05 13:41:20.813463",
    "2021-10-05 13:41:20.813466",
    "2021-10-05 13:41:20.813468",
    "2021-10-05 13:41:20.813471",
    "2021-10-05 13:41:20.813473",
    "2021-10-05 13:41:20.813476",
    "2021-10-05 13:41:20.813478",
    "2021-10-05 13:41:20.813482",
    "2021-10-05 13:41:20.813485",
    "2021-10-05 13:41:20.813488",
    "2021-10-05 13:41:20.813491",
    "2021-10-05 13:41:20.813493",
    "2021-10-05 13:41:20.813495",
    "2021-10-05 13:41:20.813498",
    "2021-10-05 13:41:20.813501",
    "2021-10-05 13:41:20.813503",
    "2021-10-05 13:41:20.813506",
    "2021-10-05 13:41:20.813508",
    "2021-10-05 13:41:20.813511",
    "2021-10-05 13:41:20.813514",
    "2021-10-05 13:41:20.813516",
    "2021-10-05 13:41:20.813519",
    "2021-10-05 13:41:20.813522",
    "2021-10-05 13:41:20.813524",
    "2021-10-05 13:41:20.813527",
    "2021-10-05 13:41:20.813529",
    "2021-10-05 13:41:20.813532",
    "2021-10-05 13:41:20.813535",
    "2021-10-05 13:41:20.813537",
    "2021-10-05 13:41:20.813540",
    "2021-10-05 13:41:20.813542",
    "2021-10-05 13:41:20.813545",
    "2021-10-05 13:41:20.813547",
    "2021-10-05 13:41:20.813550",
    "2021-10-05 13:41:20.813553",
    "2021-10-05 13:41:20.813555",
    "2021-10-05 13:41:20.813558",
    "2021-10-05 13:41:20.813560",
    "2021-10-05 13:41:20.813563",
    "2021-10-05 13:41:20.813566",
    "2021-10-05 13:41:20.813568",
    "2021-10-05 13:41:20.813571",
    "2021-10-05 13:41:20.813573",
    "2021-10-05 13:41:20.813576",
    "2021-10-05 13:41:20.813578",
    "2021-10-05 13:41:20.813581",
    "2021-10-05 13:41:20.813584",
    "2021-10-05 13:41:20.813587",
    "2021-10-05 13:41:20.813589",
    "2021-10-05 13:41:20.813592",
    "2021-10-05 13:41:20.813594",
    "2021-10-05 13:41:20.813599",
    "2021-10-05 13:41:20.813602",
    "2021-10-05 13:41:20.813605",
    "2021-10-05 13:41:20.813608",
    "2021-10-05 13:41:20.813610",
    "2021-10-05 13:41:20.813613",
    "2021-10-05 13:41:20.813616",
    "2021-10-05 13:41:20.813618",
    "2021-10-05 13:41:20.813621",
    "2021-10-05 13:41:20.813623",
    "2021-10-05 13:41:20.813626",
    "2021-10-05 13:41:20.813629",
    "2021-10-05 13:41:20.813631",
    "2021-10-05 13:41:20.813634",
    "2021-10-05 13:41:20.813637",
    "2021-10-05 13:41:20.813639",
    "2021-10-05 13:41:20.813642",
    "2021-10-05 13:41:20.813644"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

